$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update row labels (column A) - order matters for shared-string table layout ---
$ws.Range("A2").Value = "DoS Split 1"
$ws.Range("A3").Value = "DoS Split 2"
$ws.Range("A4").Value = "Probe Split 1"
$ws.Range("A5").Value = "Probe Split 2"
$ws.Range("A6").Value = "U2R Split 1"
$ws.Range("A8").Value = "R2L Split 1"
$ws.Range("A7").Value = "U2R Split 2"
$ws.Range("A9").Value = "R2L Split 2"

# --- Update header E1 text ---
$ws.Range("E1").Value = "Percent Change"

# --- Updated F1 scores (column C) for a couple of rows ---
$ws.Range("C2").Value = 0.27579990986931002
$ws.Range("C3").Value = 0.21095178519872901
$ws.Range("C6").Value = 0.25427135678391899

# --- New "Percent Change" column (E2:E9) = D / B ---
$ws.Range("E2").Formula = "=D2/B2"
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=D$r/B$r"
}

# Apply the built-in "Percent" cell style to the new column
$ws.Range("E2:E9").Style = "Percent"

# --- Column width adjustments ---
$ws.Columns.Item(1).ColumnWidth = 12.17
$ws.Columns.Item(5).ColumnWidth = 13.6

# --- Update selection to match saved view state ---
$ws.Range("H4").Select()
